$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5156313333333333
$ws.Range("H2").Value = 1.546894
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.946674666666667
$ws.Range("N2").Value = 11.840024
$ws.Range("O2").Value = 0.008082287850277156
$ws.Range("P2").Value = 0.008317683527585098
$ws.Range("Q2").Value = 2.035029120606222
$ws.Range("R2").Value = 18.315262085456
$ws.Range("S2").Value = 0.008082287850277156
$ws.Range("T2").Value = 0.008317683527585098

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5156313333333333
$ws.Range("H3").Value = 1.546894
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 185.8027443333333
$ws.Range("N3").Value = 557.408233
$ws.Range("O3").Value = 0.3805003933455167
$ws.Range("P3").Value = 0.3915824222792467
$ws.Range("Q3").Value = 95.8057167975891
$ws.Range("R3").Value = 862.251451178302
$ws.Range("S3").Value = 0.3805003933455167
$ws.Range("T3").Value = 0.3915824222792467

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5156313333333333
$ws.Range("H4").Value = 1.546894
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 117.3394243333333
$ws.Range("N4").Value = 352.018273
$ws.Range("O4").Value = 0.2402962199184265
$ws.Range("P4").Value = 0.2472948188906589
$ws.Range("Q4").Value = 60.50388382156245
$ws.Range("R4").Value = 544.5349543940621
$ws.Range("S4").Value = 0.2402962199184265
$ws.Range("T4").Value = 0.2472948188906589

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5156313333333333
$ws.Range("H5").Value = 1.546894
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 139.7641523333333
$ws.Range("N5").Value = 419.292457
$ws.Range("O5").Value = 0.2862192112890951
$ws.Range("P5").Value = 0.2945553119511906
$ws.Range("Q5").Value = 72.06677621983978
$ws.Range("R5").Value = 648.600985978558
$ws.Range("S5").Value = 0.2862192112890951
$ws.Range("T5").Value = 0.2945553119511906

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.5156313333333333
$ws.Range("H6").Value = 1.546894
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 41.458574
$ws.Range("N6").Value = 82.917148
$ws.Range("O6").Value = 0.0849018875966847
$ws.Range("P6").Value = 0.05824976335131885
$ws.Range("Q6").Value = 21.37733978971866
$ws.Range("R6").Value = 128.264038738312
$ws.Range("S6").Value = 0.0849018875966847
$ws.Range("T6").Value = 0.05824976335131885
